$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.05688009254515978
$ws.Range("C2").Value = 0.3078030679106084
$ws.Range("D2").Value = 0.1763561601164249
$ws.Range("E2").Value = 0.4199478064193512
$ws.Range("F2").Value = 0.420389698777502
$ws.Range("G2").Value = 49

$ws.Range("B3").Value = 0.2667693795774359
$ws.Range("C3").Value = 0.5466263928059274
$ws.Range("D3").Value = 0.7076554788746405
$ws.Range("E3").Value = 0.8412226095835992
$ws.Range("F3").Value = 0.8060705778746855
$ws.Range("G3").Value = 49

$ws.Range("B4").Value = 0.5431296193130369
$ws.Range("C4").Value = 0.7231827691434326
$ws.Range("D4").Value = 1.246365227305445
$ws.Range("E4").Value = 1.116407285584184
$ws.Range("F4").Value = 0.985706586491255
$ws.Range("G4").Value = 48

$ws.Range("B5").Value = 0.4692704090458715
$ws.Range("C5").Value = 0.7347082924362607
$ws.Range("D5").Value = 1.284144682252101
$ws.Range("E5").Value = 1.133201077590425
$ws.Range("F5").Value = 1.042621157342553
$ws.Range("G5").Value = 47

$ws.Range("B6").Value = 0.3883654118239238
$ws.Range("C6").Value = 0.6070310945326703
$ws.Range("D6").Value = 0.8839353291870329
$ws.Range("E6").Value = 0.9401783496693767
$ws.Range("F6").Value = 0.8656782987275211
$ws.Range("G6").Value = 46

$ws.Range("B7").Value = 0.3703187222185918
$ws.Range("C7").Value = 0.5607740953733055
$ws.Range("D7").Value = 0.5932288700479009
$ws.Range("E7").Value = 0.7702135223740888
$ws.Range("F7").Value = 0.685502670612526
$ws.Range("G7").Value = 34

$ws.Range("B8").Value = 0.3586222435103067
$ws.Range("C8").Value = 0.5485688517592353
$ws.Range("D8").Value = 0.6215876060107235
$ws.Range("E8").Value = 0.788408273682312
$ws.Range("F8").Value = 0.7130099896635785
$ws.Range("G8").Value = 33

$ws.Range("B9").Value = 0.1822978435501495
$ws.Range("C9").Value = 0.5070926453844893
$ws.Range("D9").Value = 0.3970354312988059
$ws.Range("E9").Value = 0.6301074759902519
$ws.Range("F9").Value = 0.6229417758010422
$ws.Range("G9").Value = 16

$ws.Range("B10").Value = 0.203264939200266
$ws.Range("C10").Value = 0.4823557971524242
$ws.Range("D10").Value = 0.379312118145266
$ws.Range("E10").Value = 0.6158832017073254
$ws.Range("F10").Value = 0.6128217818122421
$ws.Range("G10").Value = 10
